# SubRES_New_Techs.xlsx - "Add files via upload"
#
# Adds a new Wind-energy technology set, mirroring the existing Solar/PV
# entries already present in the template:
#   - SEC_Comm:      new commodity  "WIND" / "Wind energy"
#   - SEC_Processes: new process    "NEW_WIND_PP" / "New Wind Power Plant"
#                    new mining set "MIN_WIND" / "Wind energy supply"
#   - MIN_IMP:       new mining-supply row referencing the new process/commodity
#   - PP:            new power-plant row referencing the new process/commodity
#
# Cell values are written in the same order the original author's Excel
# session would have produced them (SEC_Comm first, then SEC_Processes,
# then the two formula-driven sheets) so the shared-string table grows in
# the same sequence.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SEC_Comm: row 8 - Wind commodity (mirrors row 7, the Solar commodity)
# ---------------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item("SEC_Comm")
$wsComm.Range("B8").Value = "NRG"
$wsComm.Range("C8").Value = "WIND"
$wsComm.Range("D8").Value = "Wind energy"
$wsComm.Range("E8").Value = "PJ"
$wsComm.Range("G8").Value = "SEASON"

# ---------------------------------------------------------------------------
# SEC_Processes: row 9 - New Wind Power Plant process (mirrors row 8, PV)
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("SEC_Processes")
$wsProc.Range("B9").Value = "ELE"
$wsProc.Range("C9").Value = "PL"
$wsProc.Range("D9").Value = "NEW_WIND_PP"
$wsProc.Range("E9").Value = "New Wind Power Plant"
$wsProc.Range("F9").Value = "PJ"
$wsProc.Range("G9").Value = "GW"
$wsProc.Range("H9").Value = "DAYNITE"

# SEC_Processes: row 10 - Wind mining/supply set (mirrors row 7, MIN_SOLAR)
$wsProc.Range("B10").Value = "MN"
$wsProc.Range("C10").Value = "PL"
$wsProc.Range("D10").Value = "MIN_WIND"
$wsProc.Range("E10").Value = "Wind energy supply"
$wsProc.Range("F10").Value = "PJ"
$wsProc.Range("G10").Value = "PJ/a"
$wsProc.Range("H10").Value = "SEASON"

# ---------------------------------------------------------------------------
# MIN_IMP: row 9 - Wind energy supply row (mirrors row 8, MIN_SOLAR)
# ---------------------------------------------------------------------------
$wsMin = $wb.Worksheets.Item("MIN_IMP")
$wsMin.Range("B9").Formula = "=SEC_Processes!D10"
$wsMin.Range("D9").Formula = "=SEC_Comm!C8"
$wsMin.Range("E9").Value = 0.001
$wsMin.Range("F9").Value = 0.001
$wsMin.Range("G9").Value = 2025
$wsMin.Range("H9").Value = 100

# ---------------------------------------------------------------------------
# PP: row 9 - New Wind Power Plant row (mirrors row 8, NEW_PV_PP)
# ---------------------------------------------------------------------------
$wsPP = $wb.Worksheets.Item("PP")
$wsPP.Range("B9").Formula = "=SEC_Processes!D9"
$wsPP.Range("C9").Formula = "=SEC_Processes!E9"
$wsPP.Range("D9").Formula = "=SEC_Comm!C8"
$wsPP.Range("E9").Formula = "=SEC_Comm!C27"
$wsPP.Range("F9").Value = 1
$wsPP.Range("G9").Value = 31.536000000000001
$wsPP.Range("H9").Value = 0.12
$wsPP.Range("I9").Value = 1
$wsPP.Range("K9").Value = 1000
$wsPP.Range("L9").Value = 2025
$wsPP.Range("M9").Value = 25
